# Apply cryptos list update per commit "Updated cryptos list on Sun Mar 24 19:00:02 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-coerced to a number by Excel
# (e.g. "562.13", "0.999") need an explicit text format so they are stored
# as strings, matching the original inline-string cell contents.
$textCells = @(
    "D5",
    "D6",
    "D7",
    "D12",
    "D14",
    "D20",
    "D22",
    "D24",
    "D25",
    "D26",
    "D28",
    "D30",
    "D31",
    "D32",
    "D33",
    "D34",
    "D37",
    "D39",
    "D40",
    "D48",
    "D49",
    "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "65.729.82"
$ws.Range("E2").Value = "  +0.47%  "

$ws.Range("D3").Value = "3.383.51"
$ws.Range("E3").Value = "  -0.72%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "562.13"
$ws.Range("E5").Value = "  +0.03%  "

$ws.Range("D6").Value = "175.96"
$ws.Range("E6").Value = "  -0.11%  "

$ws.Range("D7").Value = "0.630"
$ws.Range("E7").Value = "  +0.21%  "

$ws.Range("D8").Value = "3.376.26"
$ws.Range("E8").Value = "  -0.67%  "

$ws.Range("E9").Value = "  +0.04%  "

$ws.Range("E10").Value = "  +1.23%  "

$ws.Range("E11").Value = "  -0.20%  "

$ws.Range("D12").Value = "53.69"
$ws.Range("E12").Value = "  -2.37%  "

$ws.Range("E13").Value = "  -0.98%  "

$ws.Range("D14").Value = "9.20"
$ws.Range("E14").Value = "  +0.35%  "

$ws.Range("D15").Value = "3.922.36"
$ws.Range("E15").Value = "  -0.69%  "

$ws.Range("E16").Value = "  -1.05%  "

$ws.Range("E17").Value = "  -0.04%  "

$ws.Range("D18").Value = "3.382.66"
$ws.Range("E18").Value = "  -0.50%  "

$ws.Range("D19").Value = "65.547.80"
$ws.Range("E19").Value = "  +0.14%  "

$ws.Range("D20").Value = "11.84"
$ws.Range("E20").Value = "  -0.68%  "

$ws.Range("E21").Value = "  -0.12%  "

$ws.Range("D22").Value = "466.35"
$ws.Range("E22").Value = "  -1.58%  "

$ws.Range("E23").Value = "  -2.69%  "

$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").Value = "14.39"
$ws.Range("E24").Value = "  +6.71%  "

$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "89.81"
$ws.Range("E25").Value = "  +2.99%  "

$ws.Range("D26").Value = "4.09"
$ws.Range("E26").Value = "  -1.45%  "

$ws.Range("E27").Value = "  -0.08%  "

$ws.Range("D28").Value = "10.61"
$ws.Range("E28").Value = "  -3.12%  "

$ws.Range("E29").Value = "  -2.17%  "

$ws.Range("D30").Value = "31.07"
$ws.Range("E30").Value = "  -0.70%  "

$ws.Range("D31").Value = "6.57"
$ws.Range("E31").Value = "  -2.90%  "

$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").Value = "581.49"
$ws.Range("E32").Value = "  +0.98%  "

$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").Value = "11.42"
$ws.Range("E33").Value = "  -1.24%  "

$ws.Range("D34").Value = "62.01"
$ws.Range("E34").Value = "  -0.87%  "

$ws.Range("E35").Value = "  -0.68%  "

$ws.Range("E36").Value = "  -0.01%  "

$ws.Range("D37").Value = "3.60"
$ws.Range("E37").Value = "  +2.02%  "

$ws.Range("E38").Value = "  +0.96%  "

$ws.Range("D39").Value = "35.96"
$ws.Range("E39").Value = "  +0.08%  "

$ws.Range("D40").Value = "0.375"
$ws.Range("E40").Value = "  +0.25%  "

$ws.Range("D41").Value = "0.0₃0741"
$ws.Range("E41").Value = "  -2.64%  "

$ws.Range("D42").Value = "3.098.29"
$ws.Range("E42").Value = "  +0.15%  "

$ws.Range("E43").Value = "  -1.36%  "

$ws.Range("E44").Value = "  -0.65%  "

$ws.Range("E45").Value = "  -1.02%  "

$ws.Range("E46").Value = "  -2.17%  "

$ws.Range("E47").Value = "  -1.37%  "

$ws.Range("D48").Value = "0.999"
$ws.Range("E48").Value = "  -0.02%  "

$ws.Range("D49").Value = "140.05"
$ws.Range("E49").Value = "  +1.17%  "

$ws.Range("E50").Value = "  -1.81%  "

$ws.Range("D51").Value = "8.48"
$ws.Range("E51").Value = "  +1.66%  "

